$d = $word.ActiveDocument

# 1) Update the address line (split into several runs in the target,
#    but a plain Find/Replace on the whole text achieves the same
#    visible result).
$d.Content.Find.Execute(
    "2030 Hawthorn Ct. Dr., Apt. #7416C, Ames IA 50012 | 630-207-0071 | matthoskinsghs@gmail.com",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "14 Wakefield Lane, Geneva IL 50012 | 630-207-0071 | matthoskinsghs@gmail.com",
    2)

# 2) Update the date line to use [Month] [day] placeholders.
$d.Content.Find.Execute(
    "April 5, 2021",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[Month] [day], 2021",
    2)
